$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" values between row 16 (period 2410) and row 22 (period 2404)
$ws.Range("F16").Value = 91000
$ws.Range("F22").Value = 54600
